$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2.982993458012097
$ws.Range("B3").Value = 3.185281877226719
$ws.Range("B4").Value = 0.5000899561570841
$ws.Range("B5").Value = 0.2356382498269848
$ws.Range("B6").Value = 0.1248707319547105
$ws.Range("B7").Value = 1.804617069099212
$ws.Range("B8").Value = 0.9854420787266219
$ws.Range("B9").Value = 0.7234265434288588
$ws.Range("B10").Value = 1.804617069099212
$ws.Range("B11").Value = 0.8387175341099418
$ws.Range("B12").Value = 0.4528785591104308
$ws.Range("B13").Value = 0.8252109522374417
$ws.Range("B14").Value = 0.5730396440149269
$ws.Range("B15").Value = 1.374095329832081
$ws.Range("B16").Value = 1.684025186837135
$ws.Range("B17").Value = 1.072034271432329
$ws.Range("B18").Value = 0.6806268500810692
$ws.Range("B19").Value = 1.06643168202768
$ws.Range("B20").Value = 1.759134974917344
$ws.Range("B21").Value = 1.405857137272286
$ws.Range("B22").Value = 0.9140008136579252
$ws.Range("B23").Value = 1.415993033268747
$ws.Range("B24").Value = 0.9968082181388143
$ws.Range("B25").Value = 0.09864551507746874
$ws.Range("B26").Value = 1.55494819431765
$ws.Range("B27").Value = 0.6974394253448498
$ws.Range("B28").Value = 1.858052519626876
$ws.Range("B29").Value = 1.131880927038853
$ws.Range("B30").Value = 0.8843010928965391
$ws.Range("B31").Value = 0.5771152653830183
$ws.Range("B32").Value = 0.5771152653830183
$ws.Range("B33").Value = 0.6399297951078362
$ws.Range("B34").Value = 0.537443682613347
$ws.Range("B35").Value = 0.5000899561570841
$ws.Range("B36").Value = 0.3485095648908511
$ws.Range("B37").Value = 0.2134279081879926
$ws.Range("B38").Value = 0.4528785591104308
$ws.Range("B39").Value = 1.33563522521528
$ws.Range("B40").Value = 0.489449455794709
$ws.Range("B41").Value = 0.8429960974594664
$ws.Range("B42").Value = 0.8264765188794441
$ws.Range("B43").Value = 0.4947697059758965
$ws.Range("B44").Value = 0.3563042493206986
$ws.Range("B45").Value = 1.095018225019914
$ws.Range("B46").Value = 1.804617069099212
$ws.Range("B47").Value = 1.730106732928166
$ws.Range("B48").Value = 0.9104918665247099
$ws.Range("B49").Value = 0.8459799601421594
$ws.Range("B50").Value = 0.5413315018234505
$ws.Range("B51").Value = 0.7234265434288588
$ws.Range("B52").Value = 0.7234265434288588
$ws.Range("B53").Value = 1.403978974286529
$ws.Range("B54").Value = 1.845312545345013
$ws.Range("B55").Value = 1.375318573454345
$ws.Range("B56").Value = 0.7928173930677473
$ws.Range("B57").Value = 1.730106732928166
$ws.Range("B58").Value = 1.069346497970898
$ws.Range("B59").Value = 1.403978974286529
$ws.Range("B60").Value = 1.730106732928166
$ws.Range("B61").Value = 0.6264547007678573
$ws.Range("B62").Value = 3.190535134532942
$ws.Range("B63").Value = 1.299446308105855
$ws.Range("B64").Value = 0.4242997605239676
$ws.Range("B65").Value = 1.858052519626876
$ws.Range("B66").Value = 1.095018225019914
$ws.Range("B67").Value = 0.5220658354212935
$ws.Range("B68").Value = 0.1998367827839339
$ws.Range("B69").Value = 1.900160076847982
$ws.Range("B70").Value = 0.2235775399841599
$ws.Range("B71").Value = 0.1998367827839339
$ws.Range("B72").Value = 1.573281009989405
$ws.Range("B73").Value = 1.812666695838897
$ws.Range("B74").Value = 0.489449455794709
$ws.Range("B75").Value = 0.5086371355091618
$ws.Range("B76").Value = 0.4684272258084257
$ws.Range("B77").Value = 0.5674970173644726
$ws.Range("B78").Value = 1.030943476651703
$ws.Range("B79").Value = 2.063481559231944
$ws.Range("B80").Value = 0.3591002507456232
$ws.Range("B81").Value = 0.9848190030180373
$ws.Range("B82").Value = 1.019558080821817
$ws.Range("B83").Value = 0.5220658354212935
$ws.Range("B84").Value = 0.303333355177287
$ws.Range("B85").Value = 1.573281009989405
$ws.Range("B86").Value = 0.5000899561570841
$ws.Range("B87").Value = 1.405919626969103
$ws.Range("B88").Value = 1.954409497046345
$ws.Range("B89").Value = 0.7614428806699844
$ws.Range("B90").Value = 0.3999049564437642
